# Update the workbook as described by the commit diff:
#  - Column C ("Förändrad") for every data row (2..39) moves from 45202 to 45203.
#  - Row 39 picks up an explicit ht="15" customHeight="1" (matches the other rows).
#  - A brand-new data row 40 is appended (case "A 46919-2023").
#  - The sheet's used-range dimension grows from A1:Y39 to A1:Y40 automatically
#    once the new row is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" (changed) date column for all existing data rows.
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# 2) Row 39 gains an explicit custom row height of 15 (same as every other row).
$ws.Rows.Item(39).RowHeight = 15

# 3) Append the new row 40 with its values.
$newRow = 40

$ws.Cells.Item($newRow, 1).Value = "A 46919-2023"

$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 2).Value = 45201

$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 3).Value = 45203

$ws.Cells.Item($newRow, 4).Value = "STOCKHOLMS LÄN"
$ws.Cells.Item($newRow, 5).Value = "SALEM"
$ws.Cells.Item($newRow, 6).Value = "Kommuner"

$ws.Cells.Item($newRow, 7).Value = 1.3
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# R40 stays blank but, like the rest of column R, wraps text.
$ws.Cells.Item($newRow, 18).WrapText = $true

Write-Output "Applied Salem update: C2:C39 -> 45203, row39 height, row40 appended."
